# Apply "add api for message" changes to the "message" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")

# --- Content edits ---
# A4: clarify that "sys" (not "all") denotes the system log type
$ws.Range("A4").Value = "sys表示是系统日志"

# B4: new column documenting from_id as a numeric user id
$ws.Range("B4").Value = "数字,表示用户id，"

# B5: clarify meaning of from_id == 0 (message sent by the system)
$ws.Range("B5").Value = "为0时表示是系统发出的"

# F5: clarify that "time" includes both date and time
$ws.Range("F5").Value = "包括日期和时间"

# --- Column B width grew to fit the new descriptive text ---
$ws.Columns("B").ColumnWidth = 15.75

# --- Selection ends up on F5 ---
$ws.Range("F5").Select() | Out-Null
